# edit.ps1 — applies the changes described by the commit:
#   * Fix #34: Remove OnToggle and OffToggle because is not necessary.
#   * Fix issue with Diagnosticlistener on non asp.net core apps
#       (shows up in the deck as: "deployment" -> "delivery",
#        "practice" -> "strategy" on the vision slide)
#   * Footer date field bumped from 18/09/2019 to 20/09/2019 across the
#     slide master and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the "datetimeFigureOut" footer field text: 18/09/2019 -> 20/09/2019
#    on the slide master and on every custom (slide) layout.
# ---------------------------------------------------------------------
function Update-DateField($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText -and $tf.TextRange.Text -eq "18/09/2019") {
                $tf.TextRange.Text = "20/09/2019"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateField $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 22 ("Esquio out of box toggles"): drop the "OnToggle" and
#    "OffToggle" paragraphs from the content placeholder, leaving
#    "FromToToggle" as the new first line.
# ---------------------------------------------------------------------
$slide22 = $p.Slides.Item(22)
foreach ($shp in $slide22.Shapes) {
    if ($shp.Name -eq "Content Placeholder 2") {
        $tr = $shp.TextFrame.TextRange
        # Paragraph .Text carries a trailing CR (`r) for every paragraph but
        # the very last one, so trim it before comparing literal text.
        while ($tr.Paragraphs().Count -gt 0 -and $tr.Paragraphs(1).Text.TrimEnd("`r") -eq "OnToggle") {
            $tr.Paragraphs(1).Delete()
        }
        while ($tr.Paragraphs().Count -gt 0 -and $tr.Paragraphs(1).Text.TrimEnd("`r") -eq "OffToggle") {
            $tr.Paragraphs(1).Delete()
        }
    }
}

# ---------------------------------------------------------------------
# 3) Slide 7 ("About Feature Flags"): reword the vision bullet —
#    "... is a deployment practice." -> "... is a delivery strategy."
# ---------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
foreach ($shp in $slide7.Shapes) {
    if ($shp.Name -eq "Content Placeholder 2") {
        $tr = $shp.TextFrame.TextRange
        $para = $tr.Paragraphs(2)

        $t = $para.Text
        $idx = $t.IndexOf("deployment")
        if ($idx -ge 0) {
            $para.Characters($idx + 1, 10).Text = "delivery"
        }

        $t = $para.Text
        $idx = $t.IndexOf("practice")
        if ($idx -ge 0) {
            $para.Characters($idx + 1, 8).Text = "strategy"
        }
    }
}
